# Applies the "edited the testing doc" commit to Testing.docx.
#
# Summary of content changes (paragraph numbers refer to the 2.0..9.0
# numbered items in the body):
#   2.0 - marker placement sentence reworded/extended
#   3.0 - accident time-out sentence reworded/extended
#   4.0 - pothole sentence reworded/extended; _GoBack bookmark relocated here
#   7.0 - password changed from "admin123" to "12345" (6.0 username is
#         untouched even though it also contains "admin123")
#   9.0 - cosmetic run merge only (text unchanged)
#   trailing bookmark paragraph - the _GoBack bookmark that used to sit in
#         the last empty paragraph is removed from there (it now lives on
#         paragraph 4.0 instead)

$d = $word.ActiveDocument

# Left/right curly quote characters used throughout the document.
$lq = [char]8220
$rq = [char]8221

# --- 2.0: marker placement -------------------------------------------------
$r = $d.Content
$r.Find.Execute(
    "To place a marker, choose either a pothole or an accident, and place it in the desired area on the map.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "To place a marker, choose either a pothole or an accident, which will appear at your location and can then be dragged to the desired area on the map.",
    2) | Out-Null

# --- 3.0: accident marker time-out -----------------------------------------
$r = $d.Content
$r.Find.Execute(
    "Accidents will time out after an hour.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "An accident marker will not be visible on the map after an hour has passed.",
    2) | Out-Null

# --- 4.0: pothole marker behaviour ------------------------------------------
$r = $d.Content
$r.Find.Execute(
    "Potholes will last until the admin logs in to change the pothole to fixed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Potholes will always remain on the map, however admin can decide if they want to put the pothole as " + $lq + "fixed" + $rq + " which will turn the colour of the marker to green.",
    2) | Out-Null

# Relocate the _GoBack bookmark onto the end of the 4.0 paragraph (it used to
# live in the trailing empty paragraph at the end of the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$p4 = $d.Paragraphs(6)
$p4Text = $p4.Range.Text
$p4EndPos = $p4.Range.Start + $p4Text.Length - 1
$bmRange = $d.Range($p4EndPos - 1, $p4EndPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- 7.0: password text (6.0 username keeps "admin123") --------------------
$r = $d.Paragraphs(9).Range
$r.Find.Execute(
    "admin123",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "12345",
    2) | Out-Null

# --- 9.0: merge the split "9" / ".0" runs into a single "9.0" run ----------
$r = $d.Paragraphs(11).Range
$r.Find.Execute("9.0", $true, $false, $false, $false, $false, $true, 1, $false, "9.0", 2) | Out-Null

Write-Output "edits applied"
